$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for the Macroferia Regional
# de Talca - Poroto granado series. Insert a row at 96 (shifting every
# subsequent row down by one, through the former row 124 which becomes
# row 125) and populate it with the new data point.
$ws.Rows.Item(96).Insert()

$ws.Cells.Item(96, 1).Value = 5
$ws.Cells.Item(96, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(96, 3).Value = "Maule"
$ws.Cells.Item(96, 4).Value = 44627
$ws.Cells.Item(96, 5).Value = 7
$ws.Cells.Item(96, 6).Value = 100112030
$ws.Cells.Item(96, 7).Value = "Poroto granado"
$ws.Cells.Item(96, 8).Value = "Sin especificar"
$ws.Cells.Item(96, 9).Value = "Primera"
$ws.Cells.Item(96, 10).Value = 300
$ws.Cells.Item(96, 11).Value = 25000
$ws.Cells.Item(96, 12).Value = 25000
$ws.Cells.Item(96, 13).Value = 25000
$ws.Cells.Item(96, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(96, 15).Value = "Región del Maule"
$ws.Cells.Item(96, 16).Value = 1000
$ws.Cells.Item(96, 17).Value = 25
$ws.Cells.Item(96, 18).Value = "Hortaliza"

$ws.Cells.Item(96, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
